$wb = $excel.ActiveWorkbook
$ts = $wb.Worksheets.Item("TestSteps")

# --- Insert a new "Page Object" column before column D (shifts old D/Action_Keyword to E) ---
$ts.Columns("D").Insert()

# New column D header + width
$ts.Range("D1").Value = "Page Object "
$ts.Columns("D").ColumnWidth = 20.625

# Page Object values for the rows that have a related page object
$ts.Range("D4").Value = "btn_MyAccount"
$ts.Range("D5").Value = "txtbx_UserName"
$ts.Range("D6").Value = "txtbx_Password"
$ts.Range("D7").Value = "btn_LogIn"
$ts.Range("D9").Value = "btn_LogOut"

# Selection state on TestSteps after editing
$ts.Range("D7").Select()

# --- Add the "Test Cases" sheet after "TestSteps" ---
$tc = $wb.Worksheets.Add($null, $ts)
$tc.Name = "Test Cases"

$tc.Columns("A").ColumnWidth = 12.5
$tc.Columns("B").ColumnWidth = 55.625
$tc.Columns("C").ColumnWidth = 11.5

$tc.Range("A1").Value = "Test Case ID"
$tc.Range("B1").Value = "Description"
$tc.Range("C1").Value = "Runmode"

$tc.Range("A2").Value = "Login_01"
$tc.Range("B2").Value = "login in the online Store "
$tc.Range("C2").Value = "Yes"

$tc.Range("A3").Value = "Login_02"
$tc.Range("B3").Value = "login in the online Store "
$tc.Range("C3").Value = "No"

$tc.Range("C3").Select()
$tc.Activate()
